$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AD1:AF1 with the same style as the existing header row (AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting (bold, centered, bordered) from the last existing header cell
# onto the three new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
